# Auto-generated script applying the Tiamat_Profits.xlsx data update.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) on several
# rows across all 8 job sheets, per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 1062.5
$ws.Range("J52").Value = 1133.3334
$ws.Range("L52").Value = 3400.0002
$ws.Range("N52").Value = -3720.0002

$ws.Range("H98").Value = 1400
$ws.Range("I98").Value = 1400
$ws.Range("K98").Value = 1400
$ws.Range("M98").Value = 98

$ws.Range("H122").Value = 1400
$ws.Range("I122").Value = 1400
$ws.Range("K122").Value = 4200
$ws.Range("M122").Value = -1750


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1194.2941
$ws.Range("I122").Value = 1042.75
$ws.Range("J122").Value = 1558
$ws.Range("K122").Value = 3128.25
$ws.Range("L122").Value = 4674
$ws.Range("M122").Value = -678.25
$ws.Range("N122").Value = -9574

$ws.Range("H132").Value = 18962.111
$ws.Range("I132").Value = 22145
$ws.Range("J132").Value = 3915.7273
$ws.Range("K132").Value = 66435
$ws.Range("L132").Value = 11747.1819
$ws.Range("M132").Value = -63905
$ws.Range("N132").Value = -16807.1819


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9389.688
$ws.Range("I31").Value = 9439.16
$ws.Range("J31").Value = 9327.85
$ws.Range("K31").Value = 9439.16
$ws.Range("L31").Value = 9327.85
$ws.Range("M31").Value = -9144.16
$ws.Range("N31").Value = -9917.85

$ws.Range("H34").Value = 9389.688
$ws.Range("I34").Value = 9439.16
$ws.Range("J34").Value = 9327.85
$ws.Range("K34").Value = 9439.16
$ws.Range("L34").Value = 9327.85
$ws.Range("M34").Value = -9237.16
$ws.Range("N34").Value = -9731.85

$ws.Range("H50").Value = 10453.143
$ws.Range("J50").Value = 10453.143
$ws.Range("L50").Value = 10453.143
$ws.Range("N50").Value = -11703.143

$ws.Range("H51").Value = 9076.9
$ws.Range("J51").Value = 9529.888999999999
$ws.Range("L51").Value = 9529.888999999999
$ws.Range("N51").Value = -11001.889

$ws.Range("H58").Value = 1347.4419
$ws.Range("I58").Value = 943.04
$ws.Range("J58").Value = 1909.1111
$ws.Range("K58").Value = 943.04
$ws.Range("L58").Value = 1909.1111
$ws.Range("M58").Value = -740.04
$ws.Range("N58").Value = -2315.1111

$ws.Range("H60").Value = 8036.067
$ws.Range("J60").Value = 9286.182000000001
$ws.Range("L60").Value = 9286.182000000001
$ws.Range("N60").Value = -10308.182

$ws.Range("H61").Value = 9076.9
$ws.Range("J61").Value = 9529.888999999999
$ws.Range("L61").Value = 9529.888999999999
$ws.Range("N61").Value = -10225.889

$ws.Range("H68").Value = 16399.666
$ws.Range("I68").Value = 10000
$ws.Range("J68").Value = 17679.6
$ws.Range("K68").Value = 10000
$ws.Range("L68").Value = 17679.6
$ws.Range("M68").Value = -9251
$ws.Range("N68").Value = -19177.6

$ws.Range("H71").Value = 16399.666
$ws.Range("I71").Value = 10000
$ws.Range("J71").Value = 17679.6
$ws.Range("K71").Value = 30000
$ws.Range("L71").Value = 53038.8
$ws.Range("M71").Value = -26256
$ws.Range("N71").Value = -60526.8

$ws.Range("H107").Value = 270.9375
$ws.Range("I107").Value = 240.15384
$ws.Range("J107").Value = 404.33334
$ws.Range("K107").Value = 240.15384
$ws.Range("L107").Value = 404.33334
$ws.Range("M107").Value = 1679.84616
$ws.Range("N107").Value = -4244.33334

$ws.Range("H122").Value = 1218.5883
$ws.Range("I122").Value = 1437.75
$ws.Range("J122").Value = 1023.7778
$ws.Range("K122").Value = 4313.25
$ws.Range("L122").Value = 3071.3334
$ws.Range("M122").Value = -1863.25
$ws.Range("N122").Value = -7971.3334

$ws.Range("H132").Value = 23515.979
$ws.Range("I132").Value = 29679
$ws.Range("J132").Value = 1945.4
$ws.Range("K132").Value = 89037
$ws.Range("L132").Value = 5836.200000000001
$ws.Range("M132").Value = -86507
$ws.Range("N132").Value = -10896.2

$ws.Range("H136").Value = 1347.4419
$ws.Range("I136").Value = 943.04
$ws.Range("J136").Value = 1909.1111
$ws.Range("K136").Value = 2829.12
$ws.Range("L136").Value = 5727.3333
$ws.Range("M136").Value = -279.1199999999999
$ws.Range("N136").Value = -10827.3333


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 90.26087
$ws.Range("J12").Value = 26.235294
$ws.Range("L12").Value = 78.705882
$ws.Range("N12").Value = -424.705882

$ws.Range("H131").Value = 756.1900000000001
$ws.Range("I131").Value = 482.5
$ws.Range("J131").Value = 767.59375
$ws.Range("K131").Value = 1447.5
$ws.Range("L131").Value = 2302.78125
$ws.Range("M131").Value = 3592.5
$ws.Range("N131").Value = -12382.78125


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1173.25
$ws.Range("I122").Value = 1086.6666
$ws.Range("J122").Value = 1433
$ws.Range("K122").Value = 3259.9998
$ws.Range("L122").Value = 4299
$ws.Range("M122").Value = -809.9998000000001
$ws.Range("N122").Value = -9199

$ws.Range("H132").Value = 24195.705
$ws.Range("I132").Value = 858.9655
$ws.Range("K132").Value = 2576.8965
$ws.Range("M132").Value = -46.89649999999983


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2277.6
$ws.Range("I7").Value = 2043.3334
$ws.Range("J7").Value = 2629
$ws.Range("K7").Value = 2043.3334
$ws.Range("L7").Value = 2629
$ws.Range("M7").Value = -1931.3334
$ws.Range("N7").Value = -2853

$ws.Range("H40").Value = 46880
$ws.Range("I40").Value = 1426.7646
$ws.Range("K40").Value = 1426.7646
$ws.Range("M40").Value = -1290.7646

$ws.Range("H100").Value = 49284.668
$ws.Range("I100").Value = 57132.11
$ws.Range("J100").Value = 2200
$ws.Range("K100").Value = 57132.11
$ws.Range("L100").Value = 2200
$ws.Range("M100").Value = -56591.11
$ws.Range("N100").Value = -3282

$ws.Range("H126").Value = 2277.6
$ws.Range("I126").Value = 2043.3334
$ws.Range("J126").Value = 2629
$ws.Range("K126").Value = 6130.0002
$ws.Range("L126").Value = 7887
$ws.Range("M126").Value = -3660.0002
$ws.Range("N126").Value = -12827

$ws.Range("H132").Value = 41808.04
$ws.Range("I132").Value = 53392.64
$ws.Range("J132").Value = 4158.0835
$ws.Range("K132").Value = 160177.92
$ws.Range("L132").Value = 12474.2505
$ws.Range("M132").Value = -157647.92
$ws.Range("N132").Value = -17534.2505


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5128.9614
$ws.Range("I122").Value = 2115.2144
$ws.Range("J122").Value = 8645
$ws.Range("K122").Value = 6345.6432
$ws.Range("L122").Value = 25935
$ws.Range("M122").Value = -3895.6432
$ws.Range("N122").Value = -30835

$ws.Range("H126").Value = 1033.909
$ws.Range("I126").Value = 946
$ws.Range("K126").Value = 2838
$ws.Range("M126").Value = -368

$ws.Range("H132").Value = 3957.9707
$ws.Range("I132").Value = 746.96295
$ws.Range("K132").Value = 2240.88885
$ws.Range("M132").Value = 289.1111500000002

